$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2201.5881
$ws.Range("J40").Value = 2299.5
$ws.Range("L40").Value = 2299.5
$ws.Range("N40").Value = -2649.5
$ws.Range("H98").Value = 3872
$ws.Range("J98").Value = 4998
$ws.Range("L98").Value = 4998
$ws.Range("N98").Value = -7994
$ws.Range("H113").Value = 6752.643
$ws.Range("I113").Value = 6078
$ws.Range("K113").Value = 6078
$ws.Range("M113").Value = -2824
$ws.Range("H122").Value = 3872
$ws.Range("J122").Value = 4998
$ws.Range("L122").Value = 14994
$ws.Range("N122").Value = -19894
$ws.Range("H135").Value = 3250.0715
$ws.Range("I135").Value = 3385.2307
$ws.Range("J135").Value = 1493
$ws.Range("K135").Value = 30467.0763
$ws.Range("L135").Value = 13437
$ws.Range("M135").Value = -27932.0763
$ws.Range("N135").Value = -18507
$ws.Range("H137").Value = 20565.334
$ws.Range("I137").Value = 16399.25
$ws.Range("K137").Value = 49197.75
$ws.Range("M137").Value = -46647.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 268
$ws.Range("I4").Value = 221.6
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 221.6
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -105.6
$ws.Range("N4").Value = -732
$ws.Range("H5").Value = 774.75
$ws.Range("H32").Value = 4144.827
$ws.Range("I32").Value = 4144.827
$ws.Range("K32").Value = 4144.827
$ws.Range("M32").Value = -3857.827
$ws.Range("H37").Value = 27185.2
$ws.Range("J37").Value = 28216.889
$ws.Range("L37").Value = 28216.889
$ws.Range("N37").Value = -28762.889
$ws.Range("H44").Value = 60708.855
$ws.Range("J44").Value = 60708.855
$ws.Range("L44").Value = 60708.855
$ws.Range("N44").Value = -61684.855
$ws.Range("H55").Value = 74439.11
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 74439.11
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 74439.11
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -75069.11
$ws.Range("H63").Value = 4501.2
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814
$ws.Range("H66").Value = 4501.2
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068
$ws.Range("H80").Value = 72216.89
$ws.Range("J80").Value = 72216.89
$ws.Range("L80").Value = 72216.89
$ws.Range("N80").Value = -74212.89
$ws.Range("H83").Value = 72216.89
$ws.Range("J83").Value = 72216.89
$ws.Range("L83").Value = 216650.67
$ws.Range("N83").Value = -226634.67
$ws.Range("H102").Value = 1615.1052
$ws.Range("I102").Value = 1615.1052
$ws.Range("K102").Value = 1615.1052
$ws.Range("M102").Value = 6.894800000000032
$ws.Range("H122").Value = 2611.8845
$ws.Range("I122").Value = 2588.2
$ws.Range("K122").Value = 7764.599999999999
$ws.Range("M122").Value = -5314.599999999999
$ws.Range("H132").Value = 3808.139
$ws.Range("I132").Value = 2398.0344
$ws.Range("K132").Value = 7194.1032
$ws.Range("M132").Value = -4664.1032
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 774.75
$ws.Range("H94").Value = 964.5714
$ws.Range("I94").Value = 964.5714
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 964.5714
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -513.5714
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value = 9265.871999999999
$ws.Range("I134").Value = 4347.5557
$ws.Range("K134").Value = 13042.6671
$ws.Range("M134").Value = -10507.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 506.5
$ws.Range("I22").Value = 270
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 270
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = -5700
$ws.Range("H31").Value = 3914.5715
$ws.Range("I31").Value = 1694.5
$ws.Range("K31").Value = 1694.5
$ws.Range("M31").Value = -1399.5
$ws.Range("H34").Value = 3914.5715
$ws.Range("I34").Value = 1694.5
$ws.Range("K34").Value = 1694.5
$ws.Range("M34").Value = -1492.5
$ws.Range("H74").Value = 50912.75
$ws.Range("J74").Value = 50912.75
$ws.Range("L74").Value = 50912.75
$ws.Range("N74").Value = -52660.75
$ws.Range("H77").Value = 50912.75
$ws.Range("J77").Value = 50912.75
$ws.Range("L77").Value = 152738.25
$ws.Range("N77").Value = -161474.25
$ws.Range("H94").Value = 1375.5714
$ws.Range("I94").Value = 1677.2
$ws.Range("J94").Value = 621.5
$ws.Range("K94").Value = 1677.2
$ws.Range("L94").Value = 621.5
$ws.Range("M94").Value = -1226.2
$ws.Range("N94").Value = -1523.5
$ws.Range("H141").Value = 279356.12
$ws.Range("J141").Value = 338457.84
$ws.Range("L141").Value = 338457.84
$ws.Range("N141").Value = -348817.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 157.69388
$ws.Range("I2").Value = 46.533333
$ws.Range("J2").Value = 333.21054
$ws.Range("K2").Value = 279.199998
$ws.Range("L2").Value = 1999.26324
$ws.Range("M2").Value = -166.199998
$ws.Range("N2").Value = -2225.26324
$ws.Range("H8").Value = 687.5
$ws.Range("I8").Value = 687.5
$ws.Range("K8").Value = 2062.5
$ws.Range("M8").Value = -1923.5
$ws.Range("H55").Value = 6391.85
$ws.Range("I55").Value = 1122.75
$ws.Range("J55").Value = 7709.125
$ws.Range("K55").Value = 3368.25
$ws.Range("L55").Value = 23127.375
$ws.Range("M55").Value = -3191.25
$ws.Range("N55").Value = -23481.375
$ws.Range("H68").Value = 2001
$ws.Range("I68").Value = 2000.3334
$ws.Range("K68").Value = 6001.0002
$ws.Range("M68").Value = -5190.0002
$ws.Range("H69").Value = 6778.8
$ws.Range("I69").Value = 5373.5
$ws.Range("J69").Value = 12400
$ws.Range("K69").Value = 16120.5
$ws.Range("L69").Value = 37200
$ws.Range("M69").Value = -15309.5
$ws.Range("N69").Value = -38822
$ws.Range("H71").Value = 2001
$ws.Range("I71").Value = 2000.3334
$ws.Range("K71").Value = 18003.0006
$ws.Range("M71").Value = -13947.0006
$ws.Range("H72").Value = 6778.8
$ws.Range("I72").Value = 5373.5
$ws.Range("J72").Value = 12400
$ws.Range("K72").Value = 48361.5
$ws.Range("L72").Value = 111600
$ws.Range("M72").Value = -44305.5
$ws.Range("N72").Value = -119712
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H114").Value = 438.33334
$ws.Range("I114").Value = 438.33334
$ws.Range("K114").Value = 1315.00002
$ws.Range("M114").Value = 1938.99998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 60.5
$ws.Range("I2").Value = 39
$ws.Range("J2").Value = 76.625
$ws.Range("K2").Value = 39
$ws.Range("L2").Value = 76.625
$ws.Range("M2").Value = 74
$ws.Range("N2").Value = -302.625
$ws.Range("H43").Value = 30918.846
$ws.Range("J43").Value = 35994.6
$ws.Range("L43").Value = 35994.6
$ws.Range("N43").Value = -36296.6
$ws.Range("H46").Value = 40945.2
$ws.Range("J46").Value = 42772.445
$ws.Range("L46").Value = 42772.445
$ws.Range("N46").Value = -43084.445
$ws.Range("H57").Value = 44895.2
$ws.Range("J57").Value = 44895.2
$ws.Range("L57").Value = 44895.2
$ws.Range("N57").Value = -46535.2
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3759.9
$ws.Range("I122").Value = 2900
$ws.Range("J122").Value = 4128.4287
$ws.Range("K122").Value = 8700
$ws.Range("L122").Value = 12385.2861
$ws.Range("M122").Value = -6250
$ws.Range("N122").Value = -17285.2861
$ws.Range("H136").Value = 5605.242
$ws.Range("I136").Value = 4771.35
$ws.Range("K136").Value = 14314.05
$ws.Range("M136").Value = -11764.05

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1994.5555
$ws.Range("I81").Value = 879.6
$ws.Range("J81").Value = 2423.3845
$ws.Range("K81").Value = 1759.2
$ws.Range("L81").Value = 4846.769
$ws.Range("M81").Value = -698.2
$ws.Range("N81").Value = -6968.769
$ws.Range("H84").Value = 1994.5555
$ws.Range("I84").Value = 879.6
$ws.Range("J84").Value = 2423.3845
$ws.Range("K84").Value = 8796
$ws.Range("L84").Value = 24233.845
$ws.Range("M84").Value = -3492
$ws.Range("N84").Value = -34841.845
$ws.Range("H113").Value = 623.4828
$ws.Range("I113").Value = 502.92856
$ws.Range("K113").Value = 1508.78568
$ws.Range("M113").Value = 661.21432
$ws.Range("H132").Value = 143596.39
$ws.Range("I132").Value = 232080.02
$ws.Range("J132").Value = 22936.879
$ws.Range("K132").Value = 696240.0599999999
$ws.Range("L132").Value = 68810.637
$ws.Range("M132").Value = -693710.0599999999
$ws.Range("N132").Value = -73870.637
$ws.Range("H136").Value = 6668914.5
$ws.Range("I136").Value = 11766950
$ws.Range("J136").Value = 2252.2307
$ws.Range("K136").Value = 35300850
$ws.Range("L136").Value = 6756.6921
$ws.Range("M136").Value = -35298300
$ws.Range("N136").Value = -11856.6921
